# "window handle and testng web update"
# Rename the single existing sheet to "login", add two new sheets
# ("user" and "ppt") with login/testing data, and update the saved
# window/selection state so that "user" ends up as the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename the original sheet to "login" and fix up its saved
#    selection (no more tabSelected, cursor parked at A1:B1).
# ---------------------------------------------------------------
$login = $wb.Worksheets.Item(1)
$login.Name = "login"
$login.Range("A1:B1").Select() | Out-Null

# ---------------------------------------------------------------
# 2. Add the "user" sheet (becomes physical sheet2.xml) with the
#    USERNAME/PASSWORD pair for raghava2706@gmail.com.
# ---------------------------------------------------------------
$user = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$user.Name = "user"

$user.Range("A1").Value = "USERNAME"
$user.Range("B1").Value = "PASSWORD"
$user.Range("A2").Value = "raghava2706@gmail.com"
$user.Range("B2").Value = "training"

$user.Hyperlinks.Add($user.Range("A2"), "mailto:raghava2706@gmail.com") | Out-Null

$user.Columns.Item(1).ColumnWidth = 9.67
$user.Columns.Item(2).ColumnWidth = 9.67

$user.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# 3. Add a throwaway sheet, then add "ppt" after it, then delete
#    the throwaway sheet. This advances the internal sheetId
#    counter so "ppt" lands on sheetId="4" (matching the target),
#    while it still ends up as the third and last physical sheet.
# ---------------------------------------------------------------
$placeholder = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$placeholder.Name = "placeholder"

$pptNew = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$pptNew.Name = "ppt"

$excel.DisplayAlerts = $false
$placeholder.Delete() | Out-Null
$excel.DisplayAlerts = $true

# Re-fetch by name: deleting "placeholder" shifts the underlying
# physical sheet files, and the old $pptNew reference can go stale.
$ppt = $wb.Worksheets.Item("ppt")

$ppt.Range("A1").Value = "USERNAME"
$ppt.Range("B1").Value = "PASSWORD"
$ppt.Range("C1").Value = "PPT-NAME"
$ppt.Range("D1").Value = "PPT-PWD"

$ppt.Range("A2").Value = "raghava2706@gmail.com"
$ppt.Range("B2").Value = "training"
# Assign D2 before C2 so the shared-string table picks up
# "whiteboxqa" (index 10) ahead of "SDLC - General" (index 11).
$ppt.Range("D2").Value = "whiteboxqa"
$ppt.Range("C2").Value = "SDLC - General"

$ppt.Hyperlinks.Add($ppt.Range("A2"), "mailto:raghava2706@gmail.com") | Out-Null

$ppt.Range("F7").Select() | Out-Null

# ---------------------------------------------------------------
# 5. Final cursor position / active tab: "user" is the selected
#    sheet (H8 highlighted) when the workbook is saved.
# ---------------------------------------------------------------
$user.Range("H8").Select() | Out-Null
$user.Activate() | Out-Null

Write-Host "done"
